$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 457 (shifts existing rows 457..495 down to 458..496)
$ws.Rows.Item(457).Insert()

# Populate the newly inserted row 457 with the new record
$ws.Cells.Item(457, 1).Value2 = 10
$ws.Cells.Item(457, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(457, 3).Value2 = "La Araucanía"
$ws.Cells.Item(457, 4).Value2 = 45106
$ws.Cells.Item(457, 5).Value2 = 9
$ws.Cells.Item(457, 6).Value2 = 100112044
$ws.Cells.Item(457, 7).Value2 = "Perejil"
$ws.Cells.Item(457, 8).Value2 = "Sin especificar"
$ws.Cells.Item(457, 9).Value2 = "Primera"
$ws.Cells.Item(457, 10).Value2 = 60
$ws.Cells.Item(457, 11).Value2 = 4000
$ws.Cells.Item(457, 12).Value2 = 4000
$ws.Cells.Item(457, 13).Value2 = 4000
$ws.Cells.Item(457, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(457, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(457, 16).Value2 = 1333
$ws.Cells.Item(457, 17).Value2 = 3
$ws.Cells.Item(457, 18).Value2 = "Hortaliza"
